$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 14094
$ws1.Range("F3").Value = 13806
$ws1.Range("F8").Value = 36
$ws1.Range("F9").Value = 73
$ws1.Range("F11").Value = 2178
$ws1.Range("F12").Value = 176
$ws1.Range("F13").Value = 121
$ws1.Range("F14").Value = 97
$ws1.Range("F15").Value = 215
$ws1.Range("F17").Value = 590
$ws1.Range("F19").Value = 504
$ws1.Range("F21").Value = 31
$ws1.Range("F23").Value = 887
$ws1.Range("F24").Value = 144
$ws1.Range("F25").Value = 71
$ws1.Range("F26").Value = 26
$ws1.Range("F29").Value = 73
$ws1.Range("F30").Value = 35

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F6").Value = 113
$ws2.Range("F8").Value = 1710
$ws2.Range("G8").Value = 480
$ws2.Range("F11").Value = 1
$ws2.Range("F15").Value = 1824

$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 233
$ws3.Range("F3").Value = 144
$ws3.Range("F4").Value = 137

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 233
$ws4.Range("F3").Value = 14094
$ws4.Range("F4").Value = 13806
$ws4.Range("F9").Value = 36
$ws4.Range("F10").Value = 73
$ws4.Range("F14").Value = 2178
$ws4.Range("F15").Value = 144
$ws4.Range("F16").Value = 176
$ws4.Range("F17").Value = 176
$ws4.Range("F18").Value = 121
$ws4.Range("F19").Value = 97
$ws4.Range("F20").Value = 215
$ws4.Range("F24").Value = 113
$ws4.Range("F25").Value = 137
$ws4.Range("F26").Value = 590
$ws4.Range("F28").Value = 504
$ws4.Range("F30").Value = 31
$ws4.Range("F32").Value = 887
$ws4.Range("F34").Value = 1710
$ws4.Range("G34").Value = 480
$ws4.Range("F37").Value = 1
$ws4.Range("F39").Value = 144
$ws4.Range("F40").Value = 71
$ws4.Range("F41").Value = 26
$ws4.Range("F46").Value = 73
$ws4.Range("F47").Value = 35
$ws4.Range("F48").Value = 1824
